$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = 'questions = [
    {
        "title": "You want to monitor and maintain a history of all API activity in your Azure account.  Which solution best achieves this?",
        "ques_type": 2,
        "options": [
            "Leverage Azure Monitor for APIs and select the metric \u201crequests.\u201d",
            "Leverage Azure Monitor for APIs and select the metric \u201ccapacity.\u201d",
            "Leverage Virtual Network Flow Logs.",
            "Leverage Azure Monitor for Application Gateway."
        ],
        "score": "Leverage Azure Monitor for APIs and select the metric \u201crequests.\u201d"
    },
    {
        "title": "You have an application that is used throughout the world. You created endpoints that can be used to connect to the applications in several Azure regions, using Azure Traffic Manager as the routing solution. You want each user to connect to the endpoint with the lowest latency from the user''s location.Which routing method should you use?",
        "ques_type": 2,
        "options": [
            "Performance",
            "Geographic",
            "Weighted",
            "Priority"
        ],
        "score": "Performance"
    },
    {
        "title": "As part of your company\u2019s digital transformation, your team is in charge of containerizing your on-premise application and migrating it to Azure. You decided it is best to leverage Kubernetes to orchestrate your newly containerized application. Which of the following is the best solution to minimize operational and financial costs?",
        "ques_type": 2,
        "options": [
            "Create and manage a Kubernetes cluster on VM instances.",
            "Leverage Azure AKS, Azure\u2019s managed Kubernetes solution.",
            "Leverage open source tools to create and manage a Kubernetes cluster on VM instances.",
            "Leverage a third-party vendor to create and manage a Kubernetes cluster in your Azure account."
        ],
        "score": "Leverage Azure AKS, Azure\u2019s managed Kubernetes solution."
    },
    {
        "title": "You want to improve cost analysis for your company\u2019s Azure subscription by tagging each Azure resource with the proper cost center. Several people in the company have permissions to create resources, and you need to ensure that they don\u2019t forget to add the proper cost center information during the creation of new resources. What can you use to enforce the application of desired tags?",
        "ques_type": 2,
        "options": [
            "Azure Blueprints",
            "Azure Policy",
            "Azure Advisor",
            "Azure Budgets"
        ],
        "score": "Azure Policy"
    }
]'

$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $text
$ws.Rows("2").Delete()
$ws.Rows("1").AutoFit()
